$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Cells.Item(13, 1).Value = 251984
$ws.Cells.Item(13, 3).Value = 25
$ws.Cells.Item(13, 4).Value = 338.327868852459
$ws.Cells.Item(13, 5).Value = "2025-06-04 07:00:00"
$ws.Cells.Item(13, 6).Value = "2025-06-04 07:25:00"
$ws.Cells.Item(13, 7).Value = "2025-06-04 07:25:00"
$ws.Cells.Item(13, 8).Value = "2025-06-04 13:03:19"
$ws.Cells.Item(13, 9).Value = 20638
$ws.Cells.Item(13, 12).Value = 3
$ws.Cells.Item(13, 14).Value = 39874
$ws.Cells.Item(13, 16).Value = 39874
$ws.Cells.Item(13, 17).Value = "2025-06-10 00:00:00"
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 19).Value = 1

# Row 14
$ws.Cells.Item(14, 1).Value = 251180
$ws.Cells.Item(14, 3).Value = 25
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = "2025-06-04 13:03:19"
$ws.Cells.Item(14, 6).Value = "2025-06-04 13:28:19"
$ws.Cells.Item(14, 7).Value = "2025-06-04 13:28:19"
$ws.Cells.Item(14, 8).Value = "2025-06-04 13:28:19"
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 12).Value = 4
$ws.Cells.Item(14, 14).Value = "39887 (esterno)"
$ws.Cells.Item(14, 16).Value = 39887
$ws.Cells.Item(14, 17).Value = "2025-05-20 00:00:00"
$ws.Cells.Item(14, 18).Value = -15.5613387978125
$ws.Cells.Item(14, 19).Value = 7

# Row 15
$ws.Cells.Item(15, 1).Value = 252282
$ws.Cells.Item(15, 3).Value = 25
$ws.Cells.Item(15, 4).Value = 44.88524590163934
$ws.Cells.Item(15, 5).Value = "2025-06-04 13:28:19"
$ws.Cells.Item(15, 6).Value = "2025-06-04 13:53:19"
$ws.Cells.Item(15, 7).Value = "2025-06-04 13:53:19"
$ws.Cells.Item(15, 8).Value = "2025-06-04 14:38:12"
$ws.Cells.Item(15, 9).Value = 2738
$ws.Cells.Item(15, 12).Value = 5
$ws.Cells.Item(15, 14).Value = 39885
$ws.Cells.Item(15, 16).Value = 39885
$ws.Cells.Item(15, 17).Value = "2025-06-09 00:00:00"
$ws.Cells.Item(15, 18).Value = 0
$ws.Cells.Item(15, 19).Value = 1

# Row 16
$ws.Cells.Item(16, 1).Value = 252084
$ws.Cells.Item(16, 3).Value = 35
$ws.Cells.Item(16, 4).Value = 641
$ws.Cells.Item(16, 5).Value = "2025-06-04 14:38:12"
$ws.Cells.Item(16, 6).Value = "2025-06-05 07:13:12"
$ws.Cells.Item(16, 7).Value = "2025-06-05 07:13:12"
$ws.Cells.Item(16, 8).Value = "2025-06-06 09:54:12"
$ws.Cells.Item(16, 9).Value = 39101
$ws.Cells.Item(16, 12).Value = 2
$ws.Cells.Item(16, 14).Value = 39885
$ws.Cells.Item(16, 16).Value = 39885
$ws.Cells.Item(16, 17).Value = "2025-06-30 00:00:00"
$ws.Cells.Item(16, 18).Value = -1.412647996354166
$ws.Cells.Item(16, 19).Value = 7
